$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 80
$ws.Range("D80").Value = 44559
$ws.Range("J80").Value = 400
$ws.Range("K80").Value = 7000
$ws.Range("L80").Value = 8000
$ws.Range("M80").Value = 7500
$ws.Range("O80").Value = "Región de Arica y Parinacota"
$ws.Range("P80").Value = 107

# Row 81
$ws.Range("D81").Value = 44559
$ws.Range("I81").Value = "Segunda"
$ws.Range("J81").Value = 300
$ws.Range("K81").Value = 5000
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = 5500
$ws.Range("N81").Value = "`$/caja 100 unidades"
$ws.Range("O81").Value = "Región de Arica y Parinacota"
$ws.Range("P81").Value = 55
$ws.Range("Q81").Value = 100

# Row 82
$ws.Range("D82").Value = 44286
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 1200
$ws.Range("K82").Value = 11000
$ws.Range("L82").Value = 12000
$ws.Range("M82").Value = 11500
$ws.Range("N82").Value = "`$/caja 70 unidades"
$ws.Range("P82").Value = 164
$ws.Range("Q82").Value = 70

# Row 83
$ws.Range("D83").Value = 44279
$ws.Range("J83").Value = 700

# Row 84
$ws.Range("D84").Value = 44279
$ws.Range("I84").Value = "Segunda"
$ws.Range("J84").Value = 500
$ws.Range("K84").Value = 8500
$ws.Range("L84").Value = 9000
$ws.Range("M84").Value = 8750
$ws.Range("N84").Value = "`$/caja 100 unidades"
$ws.Range("P84").Value = 88
$ws.Range("Q84").Value = 100

# Row 85
$ws.Range("D85").Value = 44322
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 400
$ws.Range("K85").Value = 10000
$ws.Range("L85").Value = 11000
$ws.Range("M85").Value = 10500
$ws.Range("N85").Value = "`$/caja 70 unidades"
$ws.Range("P85").Value = 150
$ws.Range("Q85").Value = 70

# Row 86
$ws.Range("D86").Value = 44391
$ws.Range("K86").Value = 14000
$ws.Range("L86").Value = 15000
$ws.Range("M86").Value = 14500
$ws.Range("N86").Value = "`$/caja 60 unidades"
$ws.Range("P86").Value = 242
$ws.Range("Q86").Value = 60

# Row 87
$ws.Range("D87").Value = 44391
$ws.Range("I87").Value = "Segunda"
$ws.Range("J87").Value = 240
$ws.Range("K87").Value = 12000
$ws.Range("L87").Value = 13000
$ws.Range("M87").Value = 12500
$ws.Range("N87").Value = "`$/caja 100 unidades"
$ws.Range("P87").Value = 125
$ws.Range("Q87").Value = 100

# Row 88
$ws.Range("D88").Value = 44510
$ws.Range("J88").Value = 400
$ws.Range("K88").Value = 5500
$ws.Range("L88").Value = 6000
$ws.Range("M88").Value = 5750
$ws.Range("N88").Value = "`$/caja 70 unidades"
$ws.Range("P88").Value = 82
$ws.Range("Q88").Value = 70

# Row 89
$ws.Range("D89").Value = 44321
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 500
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 11000
$ws.Range("M89").Value = 10500
$ws.Range("N89").Value = "`$/caja 70 unidades"
$ws.Range("P89").Value = 150
$ws.Range("Q89").Value = 70

# Row 90
$ws.Range("D90").Value = 44385
$ws.Range("J90").Value = 500
$ws.Range("K90").Value = 13000
$ws.Range("L90").Value = 14000
$ws.Range("M90").Value = 13500
$ws.Range("N90").Value = "`$/caja 60 unidades"
$ws.Range("P90").Value = 225
$ws.Range("Q90").Value = 60

# Row 91
$ws.Range("D91").Value = 44385
$ws.Range("J91").Value = 360
$ws.Range("K91").Value = 11000
$ws.Range("L91").Value = 12000
$ws.Range("M91").Value = 11500
$ws.Range("P91").Value = 115

# Row 92
$ws.Range("D92").Value = 44308
$ws.Range("J92").Value = 600
$ws.Range("K92").Value = 9500
$ws.Range("L92").Value = 10000
$ws.Range("M92").Value = 9750
$ws.Range("P92").Value = 139

# Row 93
$ws.Range("D93").Value = 44308
$ws.Range("J93").Value = 400
$ws.Range("K93").Value = 7500
$ws.Range("M93").Value = 7750
$ws.Range("P93").Value = 78

# Row 94
$ws.Range("D94").Value = 44238
$ws.Range("J94").Value = 700
$ws.Range("K94").Value = 10000
$ws.Range("L94").Value = 11000
$ws.Range("M94").Value = 10500
$ws.Range("P94").Value = 150

# Row 95
$ws.Range("D95").Value = 44238
$ws.Range("J95").Value = 600
$ws.Range("K95").Value = 7000
$ws.Range("L95").Value = 8000
$ws.Range("M95").Value = 7500
$ws.Range("P95").Value = 75

# Row 96
$ws.Range("D96").Value = 44175
$ws.Range("J96").Value = 2400

# Row 97
$ws.Range("D97").Value = 44175
$ws.Range("I97").Value = "Segunda"
$ws.Range("J97").Value = 1700
$ws.Range("K97").Value = 4500
$ws.Range("L97").Value = 5000
$ws.Range("M97").Value = 4750
$ws.Range("N97").Value = "`$/caja 100 unidades"
$ws.Range("P97").Value = 48
$ws.Range("Q97").Value = 100

# Row 98
$ws.Range("D98").Value = 44188
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 2500
$ws.Range("K98").Value = 6500
$ws.Range("L98").Value = 7000
$ws.Range("M98").Value = 6750
$ws.Range("N98").Value = "`$/caja 70 unidades"
$ws.Range("P98").Value = 96
$ws.Range("Q98").Value = 70

# Row 99
$ws.Range("D99").Value = 44258
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 11000
$ws.Range("L99").Value = 12000
$ws.Range("M99").Value = 11500
$ws.Range("P99").Value = 164

# Row 100
$ws.Range("D100").Value = 44258
$ws.Range("J100").Value = 500
$ws.Range("K100").Value = 8000
$ws.Range("L100").Value = 9000
$ws.Range("M100").Value = 8500
$ws.Range("P100").Value = 85

# Row 101
$ws.Range("D101").Value = 44224
$ws.Range("J101").Value = 1100
$ws.Range("K101").Value = 9500
$ws.Range("L101").Value = 10000
$ws.Range("M101").Value = 9750
$ws.Range("P101").Value = 139

# Row 102
$ws.Range("D102").Value = 44224
$ws.Range("J102").Value = 600
$ws.Range("K102").Value = 7500
$ws.Range("L102").Value = 8000
$ws.Range("M102").Value = 7750
$ws.Range("P102").Value = 78

# Row 103
$ws.Range("D103").Value = 44195
$ws.Range("J103").Value = 800
$ws.Range("K103").Value = 7500
$ws.Range("L103").Value = 8000
$ws.Range("M103").Value = 7750
$ws.Range("N103").Value = "`$/caja 70 unidades"
$ws.Range("P103").Value = 111
$ws.Range("Q103").Value = 70

# Row 104
$ws.Range("D104").Value = 44195
$ws.Range("K104").Value = 5500
$ws.Range("L104").Value = 6000
$ws.Range("M104").Value = 5750
$ws.Range("P104").Value = 58

# Row 105
$ws.Range("D105").Value = 44371
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 12500
$ws.Range("L105").Value = 13000
$ws.Range("M105").Value = 12750
$ws.Range("N105").Value = "`$/caja 60 unidades"
$ws.Range("P105").Value = 212
$ws.Range("Q105").Value = 60

# Row 106
$ws.Range("D106").Value = 44371
$ws.Range("K106").Value = 10500
$ws.Range("L106").Value = 11000
$ws.Range("M106").Value = 10750
$ws.Range("P106").Value = 108

# Row 107 (new)
$ws.Range("A107").Value = 2
$ws.Range("B107").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C107").Value = "Coquimbo"
$ws.Range("D107").Value = 44272
$ws.Range("E107").Value = 4
$ws.Range("F107").Value = 100112043
$ws.Range("G107").Value = "Pepino ensalada"
$ws.Range("H107").Value = "Sin especificar"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 10000
$ws.Range("L107").Value = 11000
$ws.Range("M107").Value = 10500
$ws.Range("N107").Value = "`$/caja 70 unidades"
$ws.Range("O107").Value = "Provincia de Limarí"
$ws.Range("P107").Value = 150
$ws.Range("Q107").Value = 70
$ws.Range("R107").Value = "Hortaliza"
$ws.Range("D107").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 108 (new)
$ws.Range("A108").Value = 2
$ws.Range("B108").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C108").Value = "Coquimbo"
$ws.Range("D108").Value = 44272
$ws.Range("E108").Value = 4
$ws.Range("F108").Value = 100112043
$ws.Range("G108").Value = "Pepino ensalada"
$ws.Range("H108").Value = "Sin especificar"
$ws.Range("I108").Value = "Segunda"
$ws.Range("J108").Value = 400
$ws.Range("K108").Value = 8000
$ws.Range("L108").Value = 9000
$ws.Range("M108").Value = 8500
$ws.Range("N108").Value = "`$/caja 100 unidades"
$ws.Range("O108").Value = "Provincia de Limarí"
$ws.Range("P108").Value = 85
$ws.Range("Q108").Value = 100
$ws.Range("R108").Value = "Hortaliza"
$ws.Range("D108").NumberFormat = "YYYY-MM-DD HH:MM:SS"
